$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RCLA_Projects")

# Row 16: GG2570080 - Hospitalito Atitlan - Mothers and Infants -> add url
$ws.Range("D16").Value = "https://rotaryatitlan.org/hospitalito-grant-2025/"

# Row 4: GG2346063 - Amigos de Santa Cruz 1 -> closed, add url
$ws.Range("D4").Value = "https://rotaryatitlan.org/visiting-our-friends-amigos-de-santa-cruz/"
$ws.Range("F4").Value = "closed"

# Row 14: GG2459764 - Amigos de Santa Cruz 2 -> add url
$ws.Range("D14").Value = "https://rotaryatitlan.org/visiting-our-friends-amigos-de-santa-cruz/"

# Row 15: GG2567164 - Amigos de Santa Cruz 3 -> add url
$ws.Range("D15").Value = "https://rotaryatitlan.org/visiting-our-friends-amigos-de-santa-cruz/"

# Row 17: GG2574529 - Reforesting Santiago -> approved
$ws.Range("F17").Value = "approved"

# Update selection to D4
$ws.Range("D4").Select()
